$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 780, shifting existing rows 780..838 down to 781..839
$ws.Rows.Item(780).Insert()

# Populate the newly inserted row 780 with the new data record
$ws.Range("A780").Value = 6
$ws.Range("B780").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C780").Value = "Metropolitana"
$ws.Range("D780").Value = 45265
$ws.Range("E780").Value = 13
$ws.Range("F780").Value = 100112052
$ws.Range("G780").Value = "Albahaca"
$ws.Range("H780").Value = "Sin especificar"
$ws.Range("I780").Value = "Primera"
$ws.Range("J780").Value = 500
$ws.Range("K780").Value = 5000
$ws.Range("L780").Value = 6000
$ws.Range("M780").Value = 5540
$ws.Range("N780").Value = '$/docena de matas'
$ws.Range("O780").Value = "Región Metropolitana"
$ws.Range("P780").Value = 923
$ws.Range("Q780").Value = 6
$ws.Range("R780").Value = "Hortaliza"
